# GMS Data Release 1
#
# Exomiser dataType sheet schema update:
#   - family_id             -> referral_id
#   - sample_id              -> platekey
#   - assembly               -> genome_build
#   - full_brothers_affected -> full_siblings_affected
#   - full_sisters_affected  -> removed (trailing row cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "referral_id"
$ws.Range("A38").Value = "full_siblings_affected"
$ws.Range("A7").Value = "platekey"
$ws.Range("A11").Value = "genome_build"

# The old full_sisters_affected row is no longer needed - clear the trailing row.
$ws.Range("A39").ClearContents()
$ws.Range("B39").ClearContents()

# Match the author's final selection/cursor position.
$ws.Range("E10").Select()
